$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Fjerde kolonne"
$ws.Range("D2").Value = "Der bringer ændringer i sin egen branch"

$ws.Range("D2").Select()
